$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-unknown time for row 2 (9:00 AM)
$ws.Range("B2").Value = 0.375
$ws.Range("B2").NumberFormat = "h:mm AM/PM"

# Fill in the previously-unknown time for row 6 (8:30 PM)
$ws.Range("B6").Value = 0.85416666666666663
$ws.Range("B6").NumberFormat = "h:mm AM/PM"

# Fill in the previously-unknown description for row 6
$ws.Range("D6").Value = "Worked to understand what we need to do to complete phase 2"

# Add new row 7 for the next day's log entry
$ws.Range("A7").Value = "2019-09-29"
$ws.Range("C7").Value = "Thomas, Siddarth, Siddharth, Shawn"

$ws.Range("D16").Select()
